$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.811.80'
$ws.Range("E2").Value = '  -1.41%  '

# Row 3
$ws.Range("D3").Value = '2.590.07'
$ws.Range("E3").Value = '  -1.24%  '

# Row 4
$ws.Range("E4").Value = '  +0.40%  '

# Row 5
$ws.Range("D5").Value = '''590.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.10%  '

# Row 6
$ws.Range("D6").Value = '''151.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.66%  '

# Row 7
$ws.Range("E7").Value = '  +0.14%  '

# Row 8
$ws.Range("D8").Value = '''0.549'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.12%  '

# Row 9
$ws.Range("D9").Value = '2.585.95'
$ws.Range("E9").Value = '  -1.33%  '

# Row 10
$ws.Range("D10").Value = '''0.121'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.39%  '

# Row 11
$ws.Range("E11").Value = '  +0.27%  '

# Row 12
$ws.Range("D12").Value = '''5.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.10%  '

# Row 13
$ws.Range("D13").Value = '''0.344'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.71%  '

# Row 14
$ws.Range("D14").Value = '''27.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.80%  '

# Row 15
$ws.Range("D15").Value = '3.074.33'
$ws.Range("E15").Value = '  -0.18%  '

# Row 16
$ws.Range("D16").Value = '''0.0000178'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.21%  '

# Row 17
$ws.Range("D17").Value = '66.705.78'
$ws.Range("E17").Value = '  -1.27%  '

# Row 18
$ws.Range("D18").Value = '2.592.18'
$ws.Range("E18").Value = '  -0.62%  '

# Row 19
$ws.Range("D19").Value = '''364.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.06%  '

# Row 20
$ws.Range("D20").Value = '''10.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.32%  '

# Row 21
$ws.Range("D21").Value = '''7.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.00%  '

# Row 22
$ws.Range("E22").Value = '  -0.51%  '

# Row 23
$ws.Range("D23").Value = '''2.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.48%  '

# Row 24
$ws.Range("D24").Value = '''0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '

# Row 25
$ws.Range("D25").Value = '''10.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.67%  '

# Row 26
$ws.Range("D26").Value = '''67.61'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.47%  '

# Row 27
$ws.Range("D27").Value = '2.730.26'
$ws.Range("E27").Value = '  -0.73%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '

# Row 29
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").Value = '''579.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.90%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0999'
$ws.Range("E30").Value = '  -4.81%  '

# Row 31
$ws.Range("E31").Value = '  -4.56%  '

# Row 32
$ws.Range("D32").Value = '''7.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.42%  '

# Row 33
$ws.Range("D33").Value = '''1.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.01%  '

# Row 34
$ws.Range("D34").Value = '''0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.33%  '

# Row 35
$ws.Range("D35").Value = '''0.122'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.29%  '

# Row 36
$ws.Range("D36").Value = '''1.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.97%  '

# Row 37
$ws.Range("D37").Value = '''4.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.51%  '

# Row 38
$ws.Range("D38").Value = '''154.08'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.04%  '

# Row 39
$ws.Range("D39").Value = '''18.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.59%  '

# Row 40
$ws.Range("D40").Value = '''0.363'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.13%  '

# Row 41
$ws.Range("D41").Value = '''5.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.27%  '

# Row 42
$ws.Range("D42").Value = '''1.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.44%  '

# Row 43
$ws.Range("D43").Value = '''2.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.73%  '

# Row 44
$ws.Range("D44").Value = '''16.78'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.53%  '

# Row 45
$ws.Range("D45").Value = '''40.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.47%  '

# Row 46
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.06%  '

# Row 47
$ws.Range("D47").Value = '''153.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.51%  '

# Row 48
$ws.Range("D48").Value = '0.0₆0292'
$ws.Range("E48").Value = '  +0.01%  '

# Row 49
$ws.Range("D49").Value = '''3.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.95%  '

# Row 50
$ws.Range("D50").Value = '''21.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.83%  '

# Row 51
$ws.Range("D51").Value = '''0.611'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.48%  '
